# Update parametric survival model estimates/SEs and the corresponding
# variance-covariance matrices for the individual-arm fits (weibull,
# lognormal, llogis, gompertz) — values refreshed after re-running the
# parametric survival + multivariate NMA stacking step.

$wb = $excel.ActiveWorkbook

# --- est/se tables -----------------------------------------------------

$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value2 = -2.61515646363081
$ws.Range("C2").Value2 = 0.320756770614934
$ws.Range("B3").Value2 = 0.213278475676705
$ws.Range("C3").Value2 = 0.172728878584829

$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value2 = 1.85768941713106
$ws.Range("C2").Value2 = 0.336969112651298
$ws.Range("B3").Value2 = -0.993295258875856
$ws.Range("C3").Value2 = 0.157840161947571

$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value2 = -1.89524984654235
$ws.Range("C2").Value2 = 0.220566185666429
$ws.Range("B3").Value2 = 1.75319588941478
$ws.Range("C3").Value2 = 0.328357094076466

$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value2 = -2.39640552258822
$ws.Range("C2").Value2 = 0.274374733260961
$ws.Range("B3").Value2 = 0.0140410826723287
$ws.Range("C3").Value2 = 0.0322726576179848

# --- variance-covariance matrices ---------------------------------------

$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value2 = 0.102884905895322
$ws.Range("B2").Value2 = -0.0348334841073863
$ws.Range("A3").Value2 = -0.0348334841073863
$ws.Range("B3").Value2 = 0.0298352654971725

$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value2 = 0.113548182881003
$ws.Range("B2").Value2 = -0.0414908144906914
$ws.Range("A3").Value2 = -0.0414908144906914
$ws.Range("B3").Value2 = 0.0249135167236354

$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value2 = 0.0486494422594377
$ws.Range("B2").Value2 = 0.00702714895855935
$ws.Range("A3").Value2 = 0.00702714895855935
$ws.Range("B3").Value2 = 0.107818381230341

$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value2 = 0.0752814942520235
$ws.Range("B2").Value2 = -0.00475935260339732
$ws.Range("A3").Value2 = -0.00475935260339732
$ws.Range("B3").Value2 = 0.00104152442972767

Write-Output "edits applied"
